$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 4169.857
$ws.Range("I74").Value = 4148.5
$ws.Range("K74").Value = 4148.5
$ws.Range("M74").Value = -3212.5
$ws.Range("H77").Value = 4169.857
$ws.Range("I77").Value = 4148.5
$ws.Range("K77").Value = 20742.5
$ws.Range("M77").Value = -16062.5
$ws.Range("H131").Value = 2414
$ws.Range("J131").Value = 3609.6365
$ws.Range("L131").Value = 10828.9095
$ws.Range("N131").Value = -20908.9095
$ws.Range("H132").Value = 1154.7142
$ws.Range("I132").Value = 1100.4412
$ws.Range("K132").Value = 3301.3236
$ws.Range("M132").Value = -771.3235999999997
$ws.Range("H137").Value = 1806.2941
$ws.Range("I137").Value = 1525.5
$ws.Range("J137").Value = 1892.6923
$ws.Range("K137").Value = 4576.5
$ws.Range("L137").Value = 5678.0769
$ws.Range("M137").Value = -2026.5
$ws.Range("N137").Value = -10778.0769
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3663.2932
$ws.Range("I32").Value = 2273.2273
$ws.Range("J32").Value = 8032.0713
$ws.Range("K32").Value = 2273.2273
$ws.Range("L32").Value = 8032.0713
$ws.Range("M32").Value = -1986.2273
$ws.Range("N32").Value = -8606.0713
$ws.Range("H61").Value = 5469
$ws.Range("I61").Value = 2500
$ws.Range("J61").Value = 7448.3335
$ws.Range("K61").Value = 2500
$ws.Range("L61").Value = 7448.3335
$ws.Range("M61").Value = -2288
$ws.Range("N61").Value = -7872.3335
$ws.Range("H74").Value = 10000
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 10000
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 10000
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -11748
$ws.Range("H77").Value = 10000
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 10000
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 50000
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -58736
$ws.Range("H102").Value = 2000
$ws.Range("I102").Value = 2000
$ws.Range("K102").Value = 2000
$ws.Range("M102").Value = -378
$ws.Range("H132").Value = 1760.238
$ws.Range("I132").Value = 1419.3158
$ws.Range("K132").Value = 4257.9474
$ws.Range("M132").Value = -1727.9474
$ws.Range("H135").Value = 35395.8
$ws.Range("J135").Value = 35395.8
$ws.Range("L135").Value = 35395.8
$ws.Range("N135").Value = -45535.8
$ws.Range("H136").Value = 5469
$ws.Range("I136").Value = 2500
$ws.Range("J136").Value = 7448.3335
$ws.Range("K136").Value = 7500
$ws.Range("L136").Value = 22345.0005
$ws.Range("M136").Value = -4950
$ws.Range("N136").Value = -27445.0005
$ws.Range("H139").Value = 51999
$ws.Range("J139").Value = 51999
$ws.Range("L139").Value = 51999
$ws.Range("N139").Value = -62279
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2333.1667
$ws.Range("I20").Value = 2399.8
$ws.Range("K20").Value = 2399.8
$ws.Range("M20").Value = -2152.8
$ws.Range("H81").Value = 59999
$ws.Range("J81").Value = 59999
$ws.Range("L81").Value = 59999
$ws.Range("N81").Value = -62121
$ws.Range("H84").Value = 59999
$ws.Range("J84").Value = 59999
$ws.Range("L84").Value = 179997
$ws.Range("N84").Value = -190605
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("N89").ClearContents()
$ws.Range("H134").Value = 8795.68
$ws.Range("I134").Value = 9352.25
$ws.Range("K134").Value = 28056.75
$ws.Range("M134").Value = -25521.75
$ws.Range("H135").Value = 34999.5
$ws.Range("J135").Value = 34999
$ws.Range("L135").Value = 34999
$ws.Range("N135").Value = -45139
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3862.6316
$ws.Range("I31").Value = 1419.1111
$ws.Range("J31").Value = 6061.8
$ws.Range("K31").Value = 1419.1111
$ws.Range("L31").Value = 6061.8
$ws.Range("M31").Value = -1124.1111
$ws.Range("N31").Value = -6651.8
$ws.Range("H34").Value = 3862.6316
$ws.Range("I34").Value = 1419.1111
$ws.Range("J34").Value = 6061.8
$ws.Range("K34").Value = 1419.1111
$ws.Range("L34").Value = 6061.8
$ws.Range("M34").Value = -1217.1111
$ws.Range("N34").Value = -6465.8
$ws.Range("H58").Value = 1518.6842
$ws.Range("I58").Value = 1285
$ws.Range("J58").Value = 1688.6364
$ws.Range("K58").Value = 1285
$ws.Range("L58").Value = 1688.6364
$ws.Range("M58").Value = -1082
$ws.Range("N58").Value = -2094.6364
$ws.Range("H132").Value = 2524.0952
$ws.Range("I132").Value = 1311.25
$ws.Range("K132").Value = 3933.75
$ws.Range("M132").Value = -1403.75
$ws.Range("H136").Value = 1518.6842
$ws.Range("I136").Value = 1285
$ws.Range("J136").Value = 1688.6364
$ws.Range("K136").Value = 3855
$ws.Range("L136").Value = 5065.9092
$ws.Range("M136").Value = -1305
$ws.Range("N136").Value = -10165.9092
$ws.Range("H138").Value = 100000
$ws.Range("J138").Value = 100000
$ws.Range("L138").Value = 100000
$ws.Range("N138").Value = -110280
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 19757.078
$ws.Range("I4").Value = 105.55556
$ws.Range("K4").Value = 316.66668
$ws.Range("M4").Value = -204.66668
$ws.Range("H5").Value = 584.38464
$ws.Range("I5").Value = 517.4545000000001
$ws.Range("K5").Value = 1552.3635
$ws.Range("M5").Value = -1440.3635
$ws.Range("H114").Value = 2256.8333
$ws.Range("I114").Value = 80
$ws.Range("J114").Value = 2692.2
$ws.Range("K114").Value = 240
$ws.Range("L114").Value = 8076.599999999999
$ws.Range("M114").Value = 3014
$ws.Range("N114").Value = -14584.6
$ws.Range("H131").Value = 10887633
$ws.Range("J131").Value = 21756.842
$ws.Range("L131").Value = 65270.526
$ws.Range("N131").Value = -75350.526
$ws.Range("H135").Value = 584.38464
$ws.Range("I135").Value = 517.4545000000001
$ws.Range("K135").Value = 4657.0905
$ws.Range("M135").Value = -2122.0905
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 4345872
$ws.Range("I7").Value = 5409091
$ws.Range("J7").Value = 2006790
$ws.Range("K7").Value = 5409091
$ws.Range("L7").Value = 2006790
$ws.Range("M7").Value = -5408979
$ws.Range("N7").Value = -2007014
$ws.Range("H8").Value = 4345872
$ws.Range("I8").Value = 5409091
$ws.Range("J8").Value = 2006790
$ws.Range("K8").Value = 5409091
$ws.Range("L8").Value = 2006790
$ws.Range("M8").Value = -5408952
$ws.Range("N8").Value = -2007068
$ws.Range("H102").Value = 1554.4138
$ws.Range("I102").Value = 1361.1765
$ws.Range("J102").Value = 1828.1666
$ws.Range("K102").Value = 1361.1765
$ws.Range("L102").Value = 1828.1666
$ws.Range("M102").Value = 260.8235
$ws.Range("N102").Value = -5072.1666
$ws.Range("H113").Value = 1233.9286
$ws.Range("I113").Value = 975.625
$ws.Range("K113").Value = 975.625
$ws.Range("M113").Value = 1194.375
$ws.Range("H132").Value = 2806.6943
$ws.Range("I132").Value = 2549.7932
$ws.Range("J132").Value = 3871
$ws.Range("K132").Value = 7649.3796
$ws.Range("L132").Value = 11613
$ws.Range("M132").Value = -5119.3796
$ws.Range("N132").Value = -16673
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 9335.1875
$ws.Range("I40").Value = 6334.2
$ws.Range("J40").Value = 10699.272
$ws.Range("K40").Value = 6334.2
$ws.Range("L40").Value = 10699.272
$ws.Range("M40").Value = -6198.2
$ws.Range("N40").Value = -10971.272
$ws.Range("H61").Value = 3162.5454
$ws.Range("I61").Value = 2754.2222
$ws.Range("K61").Value = 2754.2222
$ws.Range("M61").Value = -2552.2222
$ws.Range("H74").Value = 12525000
$ws.Range("I74").Value = 25000000
$ws.Range("J74").Value = 50000
$ws.Range("K74").Value = 25000000
$ws.Range("L74").Value = 50000
$ws.Range("M74").Value = -24999002
$ws.Range("N74").Value = -51996
$ws.Range("H77").Value = 12525000
$ws.Range("I77").Value = 25000000
$ws.Range("J77").Value = 50000
$ws.Range("K77").Value = 75000000
$ws.Range("L77").Value = 150000
$ws.Range("M77").Value = -74995008
$ws.Range("N77").Value = -159984
$ws.Range("H113").Value = 3162.5454
$ws.Range("I113").Value = 2754.2222
$ws.Range("K113").Value = 2754.2222
$ws.Range("M113").Value = -584.2222000000002
$ws.Range("H122").Value = 6270.04
$ws.Range("I122").Value = 4567.8184
$ws.Range("J122").Value = 7607.5
$ws.Range("K122").Value = 13703.4552
$ws.Range("L122").Value = 22822.5
$ws.Range("M122").Value = -11253.4552
$ws.Range("N122").Value = -27722.5
$ws.Range("H134").Value = 47904
$ws.Range("J134").Value = 47904
$ws.Range("L134").Value = 47904
$ws.Range("N134").Value = -58044
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1166.5
$ws.Range("I107").Value = 999.8
$ws.Range("K107").Value = 2999.4
$ws.Range("M107").Value = -1079.4
$ws.Range("H122").Value = 15782.972
$ws.Range("I122").Value = 21357.68
$ws.Range("K122").Value = 64073.04
$ws.Range("M122").Value = -61623.04
$ws.Range("H132").Value = 3404.3914
$ws.Range("I132").Value = 1128.8572
$ws.Range("K132").Value = 3386.5716
$ws.Range("M132").Value = -856.5715999999998

Write-Host "Applied all changes"